$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 30.16742239378572
$ws.Range("B4").Value = 0.009999990463256836
$ws.Range("B6").Value = 30.16742239378571
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 0

$ws = $wb.Worksheets.Item("x")
$ws.Range("B2").Value = 2
$ws.Range("B4").Value = 6
$ws.Range("B5").Value = 7
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 12
$ws.Range("B10").Value = 3
$ws.Range("B11").Value = 10
$ws.Range("B13").Value = 5
$ws.Range("B14").Value = 13

$ws = $wb.Worksheets.Item("U")
$ws.Range("B5").Value = 3
$ws.Range("B13").Value = 3

$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value = 34.12074022476546
$ws.Range("B4").Value = 30
$ws.Range("B5").Value = 30.34885527085025
$ws.Range("B6").Value = 30
$ws.Range("B7").Value = 36.71579249669672
$ws.Range("B8").Value = 30
$ws.Range("B9").Value = 34.72107346555759
$ws.Range("B10").Value = 32.31224998648503
$ws.Range("B11").Value = 34.76592070603971
$ws.Range("B12").Value = 36.73266487536227
$ws.Range("B13").Value = 37.90090852477161
$ws.Range("B14").Value = 32.66758337047729
$ws.Range("B15").Value = 40.83745476036248

$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 191.4200000000015
$ws.Range("C8").Value = 198.2700000000015
$ws.Range("C9").Value = 197.9850000000015
$ws.Range("C10").Value = 199.5900000000015
$ws.Range("C11").Value = 203.15
$ws.Range("C12").Value = 81.47500000000072
$ws.Range("C13").Value = 80.68000000000072
$ws.Range("C14").Value = 84.71500000000073
$ws.Range("C15").Value = 80.43500000000073
$ws.Range("C16").Value = 87.34500000000074
$ws.Range("C17").Value = 154.3
$ws.Range("C18").Value = 148.3449999999993
$ws.Range("C19").Value = 128.7049999999993
$ws.Range("C20").Value = 146.3249999999993
$ws.Range("C21").Value = 134.2149999999993
$ws.Range("C22").Value = 72.6299999999995
$ws.Range("C23").Value = 80.0549999999995
$ws.Range("C24").Value = 82.31999999999948
$ws.Range("C25").Value = 83.9549999999995
$ws.Range("C26").Value = 80.8149999999995
$ws.Range("C27").Value = 295.9199999999996
$ws.Range("C28").Value = 323.5
$ws.Range("C29").Value = 294.2649999999996
$ws.Range("C30").Value = 311.1
$ws.Range("C31").Value = 297.3649999999996
$ws.Range("C32").Value = 107.3799999999999
$ws.Range("C33").Value = 112.2399999999999
$ws.Range("C34").Value = 93.78999999999985
$ws.Range("C35").Value = 108.8349999999998
$ws.Range("C36").Value = 94.77999999999986
$ws.Range("C37").Value = 274.4950000000024
$ws.Range("C38").Value = 282.9900000000024
$ws.Range("C39").Value = 275.9600000000024
$ws.Range("C40").Value = 289.3600000000025
$ws.Range("C41").Value = 285.0050000000024
$ws.Range("C42").Value = 140.5549999999989
$ws.Range("C43").Value = 159.2149999999988
$ws.Range("C44").Value = 142.1399999999989
$ws.Range("C45").Value = 147.7249999999989
$ws.Range("C46").Value = 139.7449999999989
$ws.Range("C47").Value = 226.0399999999994
$ws.Range("C48").Value = 247.1799999999993
$ws.Range("C49").Value = 221.8549999999994
$ws.Range("C50").Value = 238.4549999999994
$ws.Range("C51").Value = 224.4749999999994
$ws.Range("C52").Value = 332.4450000000016
$ws.Range("C53").Value = 341.6700000000017
$ws.Range("C54").Value = 337.6900000000016
$ws.Range("C55").Value = 350.0150000000017
$ws.Range("C56").Value = 337.9200000000017
$ws.Range("C57").Value = 295.9199999999996
$ws.Range("C58").Value = 323.5
$ws.Range("C59").Value = 294.2649999999996
$ws.Range("C60").Value = 311.1
$ws.Range("C61").Value = 297.3649999999996
$ws.Range("C62").Value = 154.3
$ws.Range("C63").Value = 148.3449999999993
$ws.Range("C64").Value = 128.7049999999993
$ws.Range("C65").Value = 146.3249999999993
$ws.Range("C66").Value = 134.2149999999993
$ws.Range("C67").Value = 332.4450000000016
$ws.Range("C68").Value = 341.6700000000017
$ws.Range("C69").Value = 337.6900000000016
$ws.Range("C70").Value = 350.0150000000017
$ws.Range("C71").Value = 337.9200000000017

$ws = $wb.Worksheets.Item("R")
$ws.Range("C12").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("C16").Value = 0

$ws = $wb.Worksheets.Item("L")
$ws.Range("C17").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("C21").Value = 0

# Delete rows 2:6 entirely on sheets alpha, y, rho (dimension shrinks)
$wsAlpha = $wb.Worksheets.Item("alpha")
$wsAlpha.Rows("2:6").Delete()

$wsY = $wb.Worksheets.Item("y")
$wsY.Rows("2:6").Delete()

$wsRho = $wb.Worksheets.Item("rho")
$wsRho.Rows("2:6").Delete()
